$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Insert two new rows at position 13 (pushes the old rows 13-24 down to
#    15-26, carrying their styles / row heights / shared-string refs along).
# ---------------------------------------------------------------------------
$ws.Rows("13:14").Insert()

# Remove the stray empty "A" cells that Insert() produced in the two new
# rows (rows 13/14 should only have data in columns B/C).
$ws.Range("A13").Clear()
$ws.Range("A14").Clear()

# ---------------------------------------------------------------------------
# 2) Populate the two new rows (13 and 14) with the professors' names that
#    used to live elsewhere. Column widths/styles for B & C are picked up
#    from row 2, which already carries the correct "body" styles (B -> wrap
#    normal font, C -> wrap red font).
# ---------------------------------------------------------------------------
$ws.Range("B13").Value = '5840730 - Antonio Jefferson da Silva Machado'
$ws.Range("C13").Value = '5840730 - Antonio Jefferson da Silva Machado'

$ws.Range("B14").Value = '5840726 - Cristina Bormio Nunes'
$ws.Range("C14").Value = '5840726 - Cristina Bormio Nunes'

$ws.Range("B2:C2").Copy()
$ws.Range("B13:C14").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 3) Update the text that changed value/position once the sheet was
#    reshuffled (rows are now numbered 15-26).
# ---------------------------------------------------------------------------

# Objetivos / Objectives (row 10/11 unchanged positions, but Portuguese
# Objetivos text is now a real sentence instead of the teacher name).
$ws.Range("B10").Value = 'Propiciar ao aluno os conhecimentos básicos de materiais magnéticos e supercondutores visando sua aplicação em dispositivos.'
$ws.Range("C10").Value = 'Propiciar ao aluno os conhecimentos básicos de materiais magnéticos e supercondutores visando sua aplicação em dispositivos.'

$progResumidoPt = @'
Conceitos fundamentais de propriedades magnéticas da matéria. Magnetismo de elétrons. Ferromagnetismo. Materiais magnéticos e Aplicações: moles e duros. Interação de troca em óxidos e metais. Magnetismo - Fenomenologia Clássica: diamagnetismo e paramagnetismo. Magnetismo - Fenomenologia Quântica: ferromagnetismo. Anisotropia Magnética e Interação Spin-Órbita. Magnetostricção e materiais magnetostrictivos -Introdução e aplicações. Conceitos básicos de supercondutividade. Supercondutividade - Origem Quântica. Super-onda – Consequências. Interferencia quântica – SQUID. Materiais Supercondutores e Aplicações.
'@

# Row 15 : "Programa resumido:" now holds the long Portuguese summary text
# (used to hold the activation date, which was a leftover from the old
# layout).
$ws.Range("B15").Value = $progResumidoPt
$ws.Range("C15").Value = $progResumidoPt

# Row 17 : "Programa:" also repeats the same long Portuguese summary text
# (it used to hold the teacher's name).
$ws.Range("B17").Value = $progResumidoPt
$ws.Range("C17").Value = $progResumidoPt

# Row 20 : "Método:" now holds the methodology text (used to hold the
# second teacher's name).
$ws.Range("B20").Value = 'Aulas expositivas, seminários e exercícios comentados.'
$ws.Range("C20").Value = 'Aulas expositivas, seminários e exercícios comentados.'

# Row 21 : "Critério:" now holds the grading-criteria text (used to hold
# the methodology text).
$ws.Range("B21").Value = 'A nota final , antes da recuperação é dada pela média aritmética das notas das avaliações escritas e da nota do seminário apresentado, se aplicável.'
$ws.Range("C21").Value = 'A nota final , antes da recuperação é dada pela média aritmética das notas das avaliações escritas e da nota do seminário apresentado, se aplicável.'

# Row 22 : "Norma de recuperação:" now holds the make-up-exam text (used
# to hold the grading-criteria text).
$ws.Range("B22").Value = 'Aplicação de uma prova escrita dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final da primeira avaliação'
$ws.Range("C22").Value = 'Aplicação de uma prova escrita dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final da primeira avaliação'

$bibliografia = @'
JILES, D. C. Introduction to Magnetism and Magnetic Materials, CRC Press, 1998.
COEY, J. M. D. Magnetism and Magnetic Materials, Cambridge University Press, 2010.
BUSCHOW, K. H. J.; DE BOER, F. R. Physics of Magnetism and Magnetic Materials, Springer, 2003.
CULLITY, B. D.; GRAHAM, C. D. Introduction to Magnetic Materials, Wiley-IEEE Press, 2008.
POOLE, C. P. et al., Superconductivity, Academic Press, 2007.
SHEAHEN, T. P. Introduction to High-Temperature Superconductivity, Kluwer Academic, 2002.
LEE, P. J. Engineering Superconductivity, Wiley-IEEE Press, 2001.
'@

# Row 23 : "Bibliografia:" now holds the full bibliography list (used to
# hold the old make-up-exam text).
$ws.Range("B23").Value = $bibliografia
$ws.Range("C23").Value = $bibliografia
